# Generate Report for Handback
# Update timestamps / status values produced by a re-run of the handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for rows 4-5
# (cfebaf2b-... and e97ff135-... entries) advances from 08:19:18 to 08:20:29.
$wsOverview.Range("G4").Value = "2016-09-06 08:20:29"
$wsOverview.Range("G5").Value = "2016-09-06 08:20:29"

# zh-cn sheet: Priority for rows 4-5 changes from "ht" to "mt".
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback DateTime (K)
$wsZhCn.Range("H4").Value = "2016-09-06 08:20:20"
$wsZhCn.Range("H5").Value = "2016-09-06 08:20:20"
$wsZhCn.Range("K4").Value = "2016-09-06 08:20:49"
$wsZhCn.Range("K5").Value = "2016-09-06 08:20:49"

# de-de sheet: Priority for rows 4-5 changes from "ht" to "mt".
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (H) advances same as Overview's G (shared string).
$wsDeDe.Range("H4").Value = "2016-09-06 08:20:29"
$wsDeDe.Range("H5").Value = "2016-09-06 08:20:29"

# de-de sheet: Correspond Handback DateTime (K)
$wsDeDe.Range("K4").Value = "2016-09-06 08:20:57"
$wsDeDe.Range("K5").Value = "2016-09-06 08:20:57"
